# Review_150.docx -> Review_149.docx content swap
$d = $word.ActiveDocument
$VT = [char]11   # Word's "manual line break" (<w:br/>) as plain text

# 1) Heading title
$d.Content.Find.Execute(
    "Review 150: Language Modeling Is Compression, 21.09.2023", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Review 149: [Short] Ambiguity-Aware In-Context Learning with Large Language Models, 18.09.23", 2)

# 2) Paper link (bold line)
$d.Content.Find.Execute(
    "Paper: https://arxiv.org/abs/2309.10668v2", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Paper: https://arxiv.org/abs/2205.01825v1", 2)

# 3) Secondary link line (paragraph 4)
$d.Content.Find.Execute(
    "https://huggingface.co/papers/2309.10668", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "https://arxiv.org/abs/2309.07900.pdf", 2)

# 4) Insert two new blank "Normal" paragraphs right after that link paragraph (#4)
$d.Paragraphs(4).Range.InsertParagraphAfter()
$d.Paragraphs(4).Range.InsertParagraphAfter()

# 5) Former <w:br/>-only paragraph (now paragraph 7) becomes the new intro sentence
$d.Paragraphs(7).Range.Find.Execute(
    $VT, $true, $false, $false, $false, $false,
    $true, 1, $false,
    "למידת in-context: מתברר שמודלי שפה גדולים מסוגלים ללמוד מההקשר (לפעמים אפילו לא צריך הקשר שנקרא zero-shot). למשל עבור משימת ניתוח סנטימנט אתם נותנים דוגמא אחת של (משפט, סנטימנט) ואז המודל מסתדר בעצמו. ", 2)

# 6) Paragraph 8 (was "מאמר זה משך..." + <w:br/>) becomes the "אז היום" paragraph
$d.Paragraphs(8).Range.Find.Execute(
    ("מאמר זה משך את ליבי מיד כשראיתי את הכותרת. הרי יש לא מעט קשרים בין מודלי חיזוי לבין מודלי כניסה. מאוד מאוד בגדול מודל חיזוי חזק אמור להפיק את פיצ'רים המהותיים ביותר של הדאטה הנחוצים לחיזוי ואותם הפיצ'רים ניתן לנצל לדחיסת הדאטה. " + $VT),
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "אז היום ב-#shorthebrewpapereviews סוקרים מאמר המציע שכלול של למידת in-context ומציע שיטה לבחירת דוגמת הקשר הטובה ביותר לדוגמת טסט נתונה. למה זה חשוב? יש דוגמאות די קשות ויש צורך לתת דוגמת הקשר ממש ״דומה״ כדי לכוון את המודל לכיוון הנכון. למשל עבור קטגוריות קרובות (כמו joy ו-amusement) או (אכזבה ועצב) מודל שפה עלול להתבלבל ואז חשוב לתת לו דוגמא ש״מסבירה לו״ את המשימה בדרך המיטבית ביותר. אז המאמר מציע גישה אלגנטית ואינטואטיבית לבחירה של דוגמת הקשר כזו:", 2)

# 7) Paragraph 9 (was "היום ב-#shorthebrewpapereviews...") becomes the semantic-similarity step
$d.Paragraphs(9).Range.Find.Execute(
    "היום ב-#shorthebrewpapereviews סוקרים מאמר מרתק של דיפמיינד שמראה שניתן להשתמש במודלי שפה מאומנים בתור דוחסי דאטה. איך הם בכלל ניגשו לבעיה הזו הרי זה לא לגמרי טריוויאלי איך ניתן להשתמש במודלי שפה לדחיסת דאטה. הרי מה שמודלי שפה (דקורר) יודעים לעשות הוא לחזות את התפלגות הסתברויות הטוקנים בהינתן הטוקנים הקודמים בסדרה. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "בוחרים את הדוגמאות הקרובות ביותר מבחינה סמנטית לדוגמת טסט (עם מודל pretrained המחשב ייצוג הטקסט).", 2)

# 8) Paragraph 10 (was "אז התברר שקיימת שיטת דחיסה...") becomes the top-probabilities step
$d.Paragraphs(10).Range.Find.Execute(
    "אז התברר שקיימת שיטת דחיסה שמנצלת הסתברויות אלו לדחיסת דאטה. השיטה נקראת קידוד אריתמטי(arithmetic encoding). עבור פיסת דאטה נתון השיטה מתחיל מאינטרוול [0,1]. בהגעה של הטוקן הבא מחלקים את האינטרוול לפי ההסתברויות בהתלפגות הטוקן הבא(בהינתן הקודמים) ולוקחים את האינטרוול המתאים לטוקן הבא בדאטה. לאחר הגעת הטוקן האחרון לוקחים את האינטרוול שהתקבל ולוקחים ממנו מספר שניתן לייצגו במספר המינימלי של ביטים.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "מחשבים את שתי הקטגוריות בעלות הסתברויות הגבוהות ביותר (עם הפרומפט שגורם למודל להפיק הסתברות לכל קטגוריה) עבור דוגמת הטסט ובוחרים מה דוגמאות המתויגות מהשלב הקודם כאלו עם אותן הקטגוריות.", 2)

# 9) Paragraph 11 (was the long compression paragraph w/ 2 breaks) becomes the model-mistakes step
$d.Paragraphs(11).Range.Find.Execute(
    (" סדרת ביטים זו תהווה את הדחיסה של פיסת הדאטה שלנו. אז המחברים לקחו את השיטה ובדקו עד כמה ניתן לדחוס דאטהסטים שונים מ 3 דומיינים (טקסט, תמונות, אודיו) ודחסו אותם עם הגישה הזו. לאחר מכן הם השוו את התוצאה עם הדוחסים הקלאסיים כמו gzip, LZMA2 ו- PNG (יש פרטים מעניינים איך המחברים התמודדו עם אורכי הקשר שונים בין מודלי שפה לדוחסים הקלאסיים) ויצא שמודלי שפה מצליחים להביא קצב דחיסה גבוה יותר ב 3 הדומיינים האלו." + $VT + $VT + "אבל תזכרו שיש כאן קאץ' קטן. דוחסים אלו לא לוקחים בחשבון את גודל המודל שנצטרך לשמור אותו אם נרצה לפענח את הדאטה). הדוחסים הקלאסיים הם מאוד קטנים ושם זה פחות רלוונטי. זה למעשה מביא אותנו לקונספט מתמטי מורכב הנקרא סיבוכיות קולמוגורוב שמודד את ״מורכבות״ הדאטה בתור סכום של תוכנה (המודל) לדחיסת הדאטה וגודל הדאטה אחרי הדחיסה. הם מדברים על זה קצת במאמר – ממליץ להעיף מבט."),
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "מדוגמאות מהשלב הקודם בוחרים את אלו שהמודל טועה בהם ומשתמשים בהם בתור דוגמת ההקשר.", 2)

# 10) Paragraph 12 (was empty) now gets the closing sentence
$d.Paragraphs(12).Range.Text = "משתמשים בדוגמאות אלו בשביל לבצע ניתוח סנטימנט של מודל עבור דוגמאות טסט (עם פרומפט מהונדס היטב)."

# Paragraph 13 stays an empty "Normal" paragraph (unchanged).

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
